# Regenerate save_data to use K instead of Strike#: recompute column G (K)
# values for each game-log row and write them back to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 1
    7  = 2
    8  = 0
    9  = 2
    10 = 2
    11 = 3
    12 = 3
    13 = 2
    14 = 1
    15 = 4
    16 = 2
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 2
    22 = 1
    23 = 1
    24 = 2
    25 = 1
    26 = 0
    27 = 1
    28 = 2
    29 = 0
    30 = 0
    31 = 2
    32 = 1
    33 = 2
    34 = 2
    35 = 1
    36 = 1
    37 = 2
    38 = 0
    39 = 2
    40 = 1
    41 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
